$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# The Pearson Edexcel logo images (footers) are currently named "image2.png"
# and must be renamed to "image1.png". The BTEC logo images (headers) are
# currently named "image1.jpg" and must be renamed to "image2.jpg".
# InlineShape has no writable Name property in the Word object model, so we
# temporarily convert each inline picture to a floating Shape (which does
# expose Name), rename it, then convert it back to an inline picture.

for ($hIdx = 1; $hIdx -le 2; $hIdx++) {
    $hdr = $sec.Headers.Item($hIdx)
    if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -ge 1) {
        $ishp = $hdr.Range.InlineShapes.Item(1)
        $shp = $ishp.ConvertToShape()
        $shp.Name = "image2.jpg"
        $shp.ConvertToInlineShape() | Out-Null
    }
}

for ($fIdx = 1; $fIdx -le 2; $fIdx++) {
    $ftr = $sec.Footers.Item($fIdx)
    if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -ge 1) {
        $ishp = $ftr.Range.InlineShapes.Item(1)
        $shp = $ishp.ConvertToShape()
        $shp.Name = "image1.png"
        $shp.ConvertToInlineShape() | Out-Null
    }
}
